$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Sheet1")

# --- Data change -----------------------------------------------------
# Sample row's "Tags" column (N2) used a comma to separate multiple tags;
# switch it to a semicolon-separated list.
$ws.Range("N2").Value = "Rozgar Mela; Finance"

# --- View state --------------------------------------------------------
# Scroll the window so column B is the left-most visible column, and
# leave the selection on N10 (matches the saved view state of the sheet).
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N10").Select()
